# Edit script: apply Methods.docx changes
$d = $word.ActiveDocument

# Change 1: Add first-line indent (0.5in / 720 twips) to the "Methods" heading paragraph
$d.Paragraphs(1).Format.FirstLineIndent = 36

# Change 2: Rewrite the "Agglomerative clustering..." paragraph to:
#  - split the opening sentence and insert a new sentence listing the clustering methods
#  - append a new sentence about Phnom Penh being removed prior to clustering
$rng = $d.Content
$found = $rng.Find.Execute("Agglomerative clustering was conducted to create a typology", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target paragraph text"
}
$para = $rng.Paragraphs(1).Range
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="7DD18FEE" w14:textId="2ECBC277" w:rsidR="007072F5" w:rsidRPr="007072F5" w:rsidRDefault="00AB2974"><w:r><w:t>Agglomerative clustering was conducted to create a typology for provinces based on the socioeconomic variables used in the analysis above. Several agglomerative clustering approaches were assessed</w:t></w:r><w:r><w:t>. These were single linkage, complete linkage, unweighted pair-group using arithmetic averages (UPGMA), unweighted pair-group using centroids (UPGMC), Ward’s minimum variance, and flexible clustering</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>The methods were</w:t></w:r><w:r><w:t xml:space="preserve"> compared using cophenetic correlation and Gower distance metrics, and </w:t></w:r><w:r w:rsidR="004C0A38"><w:t xml:space="preserve">the </w:t></w:r><w:r><w:t xml:space="preserve">appropriate number of clusters (k) was selected </w:t></w:r><w:r w:rsidR="004C0A38"><w:t>using</w:t></w:r><w:r><w:t xml:space="preserve"> the matrix correlation statistics </w:t></w:r><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> ADDIN ZOTERO_ITEM CSL_CITATION {"citationID":"R9LUkZIu","properties":{"formattedCitation":"(Borcard et al. 2018)","plainCitation":"(Borcard et al. 2018)","noteIndex":0},"citationItems":[{"id":2706,"uris":["http://zotero.org/users/2170232/items/5WYR7AXF"],"uri":["http://zotero.org/users/2170232/items/5WYR7AXF"],"itemData":{"id":2706,"type":"chapter","abstract":"In most cases, data exploration and the computation of association matrices are preliminary steps towards deeper analyses. In this chapter you will go further by experimenting one of the large groups of analytical methods used in ecology: clustering. Practically, you will learn how to choose among various clustering methods and compute them, apply these techniques to the Doubs River data to identify groups of sites and fish species. You will also explore two methods of constrained clustering, a powerful modelling approach where the clustering process is constrained by an external data set.","collection-title":"Use R!","container-title":"Numerical Ecology with R","event-place":"Cham","ISBN":"978-3-319-71404-2","language":"en","note":"DOI: 10.1007/978-3-319-71404-2_4","page":"59-150","publisher":"Springer International Publishing","publisher-place":"Cham","source":"Springer Link","title":"Cluster Analysis","URL":"https://doi.org/10.1007/978-3-319-71404-2_4","author":[{"family":"Borcard","given":"Daniel"},{"family":"Gillet","given":"François"},{"family":"Legendre","given":"Pierre"}],"editor":[{"family":"Borcard","given":"Daniel"},{"family":"Gillet","given":"François"},{"family":"Legendre","given":"Pierre"}],"accessed":{"date-parts":[["2021",4,30]]},"issued":{"date-parts":[["2018"]]}}}],"schema":"https://github.com/citation-style-language/schema/raw/master/csl-citation.json"} </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidRPr="00AB2974"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>(Borcard et al. 2018)</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">The capital city of Phnom Penh, which is technically a </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>province in itself, was</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> removed prior to clustering because it has extreme values for many of the variables and is thus an outlier that affects the clustering. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para.InsertXML($xml)
